$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.524.06"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "1.829.16"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.69"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("E6").Value = "  +0.20%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5179"
$ws.Range("E7").Value = "  +2.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3871"
$ws.Range("E8").Value = "  -0.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08313"
$ws.Range("E9").Value = "  +8.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.124"
$ws.Range("E10").Value = "  +0.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.92"
$ws.Range("E11").Value = "  +0.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.407"
$ws.Range("E12").Value = "  +2.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.23"
$ws.Range("E13").Value = "  +0.86%  "
$ws.Range("E14").Value = "  +0.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.519"
$ws.Range("E15").Value = "  -0.57%  "
$ws.Range("D16").Value = "1.827.82"
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "94.16"
$ws.Range("E17").Value = "  +0.65%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001125"
$ws.Range("E18").Value = "  +3.70%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06652"
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.86"
$ws.Range("E20").Value = "  +0.88%  "
$ws.Range("E21").Value = "  +0.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.073"
$ws.Range("E22").Value = "  -1.49%  "
$ws.Range("D23").Value = "28.555.64"
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.48"
$ws.Range("E24").Value = "  +3.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.252"
$ws.Range("E25").Value = "  -0.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "21.17"
$ws.Range("E26").Value = "  +2.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.29"
$ws.Range("E27").Value = "  +1.51%  "
$ws.Range("D28").Value = "2.037.26"
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.424"
$ws.Range("E29").Value = "  +0.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.21"
$ws.Range("E30").Value = "  +0.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.103"
$ws.Range("E32").Value = "  -2.68%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.07608"
$ws.Range("E33").Value = "  +7.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.742"
$ws.Range("E34").Value = "  +1.20%  "
$ws.Range("E35").Value = "  +0.74%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2237"
$ws.Range("E36").Value = "  +0.58%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02379"
$ws.Range("E37").Value = "  +2.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.312"
$ws.Range("E38").Value = "  +3.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "12.05"
$ws.Range("E39").Value = "  +7.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.804"
$ws.Range("E40").Value = "  -0.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6408"
$ws.Range("E41").Value = "  +2.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.196"
$ws.Range("E42").Value = "  +0.68%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.396"
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.70"
$ws.Range("E44").Value = "  +1.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6164"
$ws.Range("E45").Value = "  +4.27%  "
$ws.Range("E46").Value = "  +2.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "128.03"
$ws.Range("E47").Value = "  +2.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.002"
$ws.Range("E48").Value = "  +1.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.208"
$ws.Range("E49").Value = "  +1.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06983"
$ws.Range("E50").Value = "  +0.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "74.44"
$ws.Range("E51").Value = "  +0.48%  "
